$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay stored as text,
# matching the original inline-string representation of the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.952.52'
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").Value = '2.363.94'
$ws.Range("E3").Value = '  +2.27%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '302.32'
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").Value = '95.86'
$ws.Range("E6").Value = '  +0.45%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.503'
$ws.Range("E8").Value = '  -0.38%  '

$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").Value = '34.02'
$ws.Range("E10").Value = '  -0.48%  '

$ws.Range("E11").Value = '  +3.73%  '

$ws.Range("E12").Value = '  +0.28%  '

$ws.Range("D13").Value = '18.42'
$ws.Range("E13").Value = '  -2.94%  '

$ws.Range("D14").Value = '6.74'

$ws.Range("D15").Value = '2.731.28'
$ws.Range("E15").Value = '  +2.25%  '

$ws.Range("D16").Value = '2.361.48'
$ws.Range("E16").Value = '  +2.50%  '

$ws.Range("D17").Value = '0.795'
$ws.Range("E17").Value = '  +0.77%  '

$ws.Range("D18").Value = '42.934.40'
$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").Value = '11.98'
$ws.Range("E19").Value = '  -2.01%  '

$ws.Range("E20").Value = '  +2.14%  '

$ws.Range("E21").Value = '  -0.60%  '

$ws.Range("D22").Value = '67.95'

$ws.Range("D23").Value = '234.95'
$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("E24").Value = '  -4.44%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").Value = '24.56'
$ws.Range("E27").Value = '  +1.07%  '

$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").Value = '9.27'
$ws.Range("E29").Value = '  +1.59%  '

$ws.Range("D30").Value = '31.58'
$ws.Range("E30").Value = '  -1.67%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").Value = '5.03'
$ws.Range("E32").Value = '  +0.58%  '

$ws.Range("D33").Value = '17.31'
$ws.Range("E33").Value = '  -1.82%  '

$ws.Range("D34").Value = '0.0719'
$ws.Range("E34").Value = '  +2.96%  '

$ws.Range("E35").Value = '  +4.04%  '

$ws.Range("E36").Value = '  +3.60%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '125.38'
$ws.Range("E38").Value = '  -24.66%  '

$ws.Range("E39").Value = '  -1.86%  '

$ws.Range("E40").Value = '  +2.54%  '

$ws.Range("E41").Value = '  -0.58%  '

$ws.Range("D42").Value = '21.45'
$ws.Range("E42").Value = '  +2.58%  '

$ws.Range("D43").Value = '1.935.15'

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("E45").Value = '  +2.29%  '

$ws.Range("D46").Value = '9.19'
$ws.Range("E46").Value = '  -9.15%  '

$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("D48").Value = '2.589.55'
$ws.Range("E48").Value = '  +1.94%  '

$ws.Range("E49").Value = '  +2.10%  '

$ws.Range("D50").Value = '1.14'
$ws.Range("E50").Value = '  +1.81%  '

$ws.Range("D51").Value = '51.60'
$ws.Range("E51").Value = '  -3.27%  '
